{"js": "// Split the run containing the legal-reference sentence so that the\n// trailing clause \"Ley 19.937 de Autoridad Sanitaria\" becomes its own\n// run, then color that new run red (FF0000) to highlight it \u2014 matching\n// the \"[Convenios] ... resolucion a\u00f1o 2023 FIX\" edit.\nconst targetText = \"Ley 19.937 de Autoridad Sanitaria\";\n\nconst results = context.document.body.search(targetText, { matchCase: true, matchWholeWord: false });\nresults.load(\"items,text\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Target text not found: \" + targetText);\n}\n\n// Color the matched range's font red. Word/Office.js splits the\n// underlying run into separate runs automatically so only the matched\n// text (and not the text preceding it) receives the new formatting.\nconst match = results.items[0];\nmatch.font.color = \"#FF0000\";\n\nawait context.sync();\n", "ps1": "# Split the run containing the legal-reference sentence so that the\n# trailing clause \"Ley 19.937 de Autoridad Sanitaria\" becomes its own\n# run, then color that new run red (RGB 255,0,0) to highlight it \u2014\n# matching the \"[Convenios] ... resolucion a\u00f1o 2023 FIX\" edit.\n$d = $word.ActiveDocument\n\n$rng = $d.Content\n$find = $rng.Find\n$find.Text = \"Ley 19.937 de Autoridad Sanitaria\"\n$find.MatchCase = $true\n$find.Execute() | Out-Null\n\nif ($find.Found) {\n    # wdColorRed = 255 (OLE_COLOR 0x0000FF -> R=FF,G=00,B=00)\n    $rng.Font.Color = 255\n}\n"}
